$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 8.406875610351562
$ws.Range("B2").Value = -12.23353862762451
$ws.Range("C2").Value = 5.69203519821167
$ws.Range("A3").Value = -6.086453437805176
$ws.Range("B3").Value = -26.34650039672852
$ws.Range("C3").Value = -14.22333812713623
$ws.Range("A4").Value = 5.984857559204102
$ws.Range("B4").Value = -12.02742099761963
$ws.Range("C4").Value = -10.5483455657959
$ws.Range("A5").Value = -9.330942153930664
$ws.Range("B5").Value = -13.30614566802978
$ws.Range("C5").Value = -3.071574211120605
$ws.Range("A6").Value = -9.185368537902832
$ws.Range("B6").Value = -8.196459770202637
$ws.Range("C6").Value = -3.727138996124268
$ws.Range("A7").Value = -11.09085464477539
$ws.Range("B7").Value = -10.81818866729736
$ws.Range("C7").Value = 13.48275184631348
$ws.Range("A8").Value = -5.257553100585938
$ws.Range("B8").Value = -7.122329711914063
$ws.Range("C8").Value = -5.412375450134277
$ws.Range("A9").Value = 5.622566699981689
$ws.Range("B9").Value = -31.63553428649902
$ws.Range("C9").Value = 26.30605506896973
$ws.Range("A10").Value = -31.68562698364257
$ws.Range("B10").Value = 8.553699493408203
$ws.Range("C10").Value = -25.76754570007324
$ws.Range("A11").Value = -25.48157119750977
$ws.Range("B11").Value = -8.927703857421875
$ws.Range("C11").Value = -22.08453178405762
$ws.Range("A12").Value = 3.429775714874268
$ws.Range("B12").Value = -31.35194206237793
$ws.Range("C12").Value = 6.574334621429443
$ws.Range("A13").Value = -11.71933746337891
$ws.Range("B13").Value = -2.056103706359864
$ws.Range("C13").Value = 0.026987075805664
$ws.Range("A14").Value = -4.002825736999512
$ws.Range("B14").Value = 1.024898052215576
$ws.Range("C14").Value = 5.117716789245605
$ws.Range("A15").Value = -2.464081764221191
$ws.Range("B15").Value = -30.21155166625977
$ws.Range("C15").Value = 11.36835384368896
$ws.Range("A16").Value = 18.32395362854004
$ws.Range("B16").Value = 2.226233005523682
$ws.Range("C16").Value = -8.703231811523438
$ws.Range("A17").Value = -23.21679878234864
$ws.Range("B17").Value = -14.6329288482666
$ws.Range("C17").Value = -21.63663291931152
$ws.Range("A18").Value = -1.97741436958313
$ws.Range("B18").Value = -62.5509033203125
$ws.Range("C18").Value = -12.70803451538086
$ws.Range("A19").Value = -14.78749465942383
$ws.Range("B19").Value = 9.399082183837891
$ws.Range("C19").Value = 0.7060952186584473
$ws.Range("A20").Value = -45.69917297363281
$ws.Range("B20").Value = -85.84348297119141
$ws.Range("C20").Value = 16.09267616271973
$ws.Range("A21").Value = 21.77816009521484
$ws.Range("B21").Value = 1.36038064956665
$ws.Range("C21").Value = -24.16020965576172
$ws.Range("A22").Value = -22.23860740661621
$ws.Range("B22").Value = -3.703294038772583
$ws.Range("C22").Value = -13.579345703125
$ws.Range("A23").Value = -28.10436058044434
$ws.Range("B23").Value = -46.36766815185547
$ws.Range("C23").Value = 13.63930892944336
$ws.Range("A24").Value = 0.2610459327697754
$ws.Range("B24").Value = 2.697724103927612
$ws.Range("C24").Value = -2.434926986694336
$ws.Range("A25").Value = 8.421352386474609
$ws.Range("B25").Value = -4.679108619689941
$ws.Range("C25").Value = 5.642993927001953
$ws.Range("A26").Value = 12.62557125091553
$ws.Range("B26").Value = -29.58051681518555
$ws.Range("C26").Value = -5.874727249145508
$ws.Range("A27").Value = 17.53614234924316
$ws.Range("B27").Value = 6.259291648864746
$ws.Range("C27").Value = -4.211094856262207
$ws.Range("A28").Value = -11.21247386932373
$ws.Range("B28").Value = -9.254820823669434
$ws.Range("C28").Value = -18.96420669555664
$ws.Range("A29").Value = 0.2897729873657226
$ws.Range("B29").Value = -15.05545616149902
$ws.Range("C29").Value = 21.06244659423828
$ws.Range("A30").Value = -14.59784603118896
$ws.Range("B30").Value = 12.85461044311523
$ws.Range("C30").Value = 1.184048175811768
$ws.Range("A31").Value = -74.86586761474609
$ws.Range("B31").Value = -84.31977081298828
$ws.Range("C31").Value = 76.46285247802734
